# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets get their "Latest Target File", "Latest Handback File" and
#    "Latest Handback DateTime" columns populated (plus a hyperlink on the target file cell)
#  - a few columns get widened so the new values / links are readable

$wb = $excel.ActiveWorkbook

$ba95Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/081cc6427214470329dfc2c373002892e2dad66f/e2e/ba95e3e9-1c73-4c8d-8e6a-5ebf48f94167.md"
$f889Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/081cc6427214470329dfc2c373002892e2dad66f/e2e/f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.md"
$ba95Name = "ba95e3e9-1c73-4c8d-8e6a-5ebf48f94167.md"
$f889Name = "f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.md"

$statusText = "Handed back: in sync with en-US"

# Hyperlink-style blue, matching the workbook's existing custom "HyperLink" look (RGB 100,149,237 = #6495ED)
$linkColor = 15570276

# ---------------------------------------------------------------------------
# Overview sheet: the per-language status columns (zh-cn, de-de) move from
# "Ready for handoff" to "Handed back: in sync with en-US" for both rows.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Row 2 (ba95e3e9 file): Latest Target File (I2), Latest Handback File (J2),
# Latest Handback DateTime (K2)
$wsZh.Range("I2").Value = $ba95Name
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $ba95Url, "", "", $ba95Name) | Out-Null
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = $linkColor
$wsZh.Range("J2").Value = "ba95e3e9-1c73-4c8d-8e6a-5ebf48f94167.17c2f8b8d0cca3272b78b70f46934589e5cc5551.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-23 18:51:16"

# Row 3 (f889c0fd file): Latest Target File (I3), Latest Handback File (J3),
# Latest Handback DateTime (K3)
$wsZh.Range("I3").Value = $f889Name
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $f889Url, "", "", $f889Name) | Out-Null
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = $linkColor
$wsZh.Range("J3").Value = "f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.404c9f065ed2f9196b162bbac9ae528769f2b083.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-23 18:51:16"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# Row 2 (ba95e3e9 file)
$wsDe.Range("I2").Value = $ba95Name
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $ba95Url, "", "", $ba95Name) | Out-Null
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = $linkColor
$wsDe.Range("J2").Value = "ba95e3e9-1c73-4c8d-8e6a-5ebf48f94167.17c2f8b8d0cca3272b78b70f46934589e5cc5551.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-23 18:51:23"

# Row 3 (f889c0fd file)
$wsDe.Range("I3").Value = $f889Name
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $f889Url, "", "", $f889Name) | Out-Null
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = $linkColor
$wsDe.Range("J3").Value = "f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.404c9f065ed2f9196b162bbac9ae528769f2b083.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-23 18:51:23"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664
